$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the value to be stored as text even when it looks like a plain
    # number (e.g. "0.015", "-0.137"), which Excel would otherwise convert
    # to a numeric cell. Using a leading apostrophe forces text entry, then
    # resetting the style back to Normal clears the "quote prefix" flag that
    # Excel stamps on the cell format as a side effect.
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

# Row 1 - column headers (strip " Diff-in-Diff" suffix)
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "U"
$ws.Range("D1").Value = '$\pi$'
$ws.Range("E1").Value = "FFR"
$ws.Range("F1").Value = "A"

# Row 2 - "C" row
$ws.Range("A2").Value = "C"
Set-TextValue $ws.Range("C2") "0.015"
$ws.Range("D2").Value = "0.075*"
Set-TextValue $ws.Range("E2") "-0.137"
Set-TextValue $ws.Range("F2") "-0.002"

# Row 3 - "U" row
$ws.Range("A3").Value = "U"
Set-TextValue $ws.Range("B3") "0.101"
$ws.Range("D3").Value = "0.268***"
$ws.Range("E3").Value = "-0.954***"
$ws.Range("F3").Value = "0.015*"

# Row 4 - "$\pi$" row
$ws.Range("A4").Value = '$\pi$'
$ws.Range("B4").Value = "2.237*"
$ws.Range("C4").Value = "1.226***"
$ws.Range("E4").Value = "1.811***"
Set-TextValue $ws.Range("F4") "-0.006"

# Row 5 - "FFR" row
$ws.Range("A5").Value = "FFR"
Set-TextValue $ws.Range("B5") "-0.404"
$ws.Range("C5").Value = "-0.429***"
$ws.Range("D5").Value = "0.178***"
Set-TextValue $ws.Range("F5") "-0.005"

# Row 6 - "A" row
$ws.Range("A6").Value = "A"
Set-TextValue $ws.Range("B6") "-6.431"
$ws.Range("C6").Value = "8.865*"
Set-TextValue $ws.Range("D6") "-0.76"
Set-TextValue $ws.Range("E6") "-6.257"

# Row 7 - "Constant" row
Set-TextValue $ws.Range("B7") "-0.681"
$ws.Range("C7").Value = "0.406*"
Set-TextValue $ws.Range("D7") "-0.036"
Set-TextValue $ws.Range("E7") "-0.45"
$ws.Range("F7").Value = "-0.021**"

# Row 8 - "r2_adj" row (plain numeric values, not shared strings)
$ws.Range("B8").Value = 0.23
$ws.Range("C8").Value = 0.68
$ws.Range("D8").Value = 0.47
$ws.Range("E8").Value = 0.6
$ws.Range("F8").Value = 0.34
